$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1795665634674923
$ws.Range("C2").Value = 0.5696594427244582
$ws.Range("J2").Value = 0.01857585139318885
$ws.Range("P2").Value = 0.151702786377709
$ws.Range("S2").Value = 0.0804953560371517
$ws.Range("B3").Value = 0.0053475935828877
$ws.Range("C3").Value = 0.03208556149732621
$ws.Range("J3").Value = 0.0160427807486631
$ws.Range("P3").Value = 0.7165775401069518
$ws.Range("S3").Value = 0.2299465240641711
$ws.Range("J4").Value = 0.04166666666666666
$ws.Range("P4").Value = 0.7083333333333334
$ws.Range("S4").Value = 0.25
$ws.Range("P5").Value = 0.75
$ws.Range("S5").Value = 0.25
$ws.Range("B6").Value = 0.06046511627906977
$ws.Range("F6").Value = 0.08837209302325581
$ws.Range("J6").Value = 0.2232558139534884
$ws.Range("O6").Value = 0.02790697674418605
$ws.Range("Q6").Value = 0.1348837209302326
$ws.Range("R6").Value = 0.07441860465116279
$ws.Range("S6").Value = 0.3906976744186046
$ws.Range("B7").Value = 0.1084905660377359
$ws.Range("D7").Value = 0.02358490566037736
$ws.Range("F7").Value = 0.05188679245283019
$ws.Range("J7").Value = 0.1462264150943396
$ws.Range("O7").Value = 0.009433962264150943
$ws.Range("Q7").Value = 0.169811320754717
$ws.Range("R7").Value = 0.0660377358490566
$ws.Range("S7").Value = 0.4245283018867925
$ws.Range("B8").Value = 0.1102204408817635
$ws.Range("D8").Value = 0.02004008016032064
$ws.Range("F8").Value = 0.06412825651302605
$ws.Range("J8").Value = 0.12625250501002
$ws.Range("O8").Value = 0.01202404809619238
$ws.Range("Q8").Value = 0.1743486973947896
$ws.Range("R8").Value = 0.09018036072144289
$ws.Range("S8").Value = 0.4028056112224449
$ws.Range("B9").Value = 0.09248554913294797
$ws.Range("D9").Value = 0.02312138728323699
$ws.Range("E9").Value = 0.005780346820809248
$ws.Range("F9").Value = 0.05780346820809248
$ws.Range("J9").Value = 0.08092485549132948
$ws.Range("O9").Value = 0.01734104046242774
$ws.Range("Q9").Value = 0.1560693641618497
$ws.Range("R9").Value = 0.1329479768786127
$ws.Range("S9").Value = 0.4335260115606936
$ws.Range("B10").Value = 0.1321070234113712
$ws.Range("D10").Value = 0.02675585284280936
$ws.Range("E10").Value = 0.002508361204013378
$ws.Range("F10").Value = 0.06605351170568562
$ws.Range("J10").Value = 0.1187290969899666
$ws.Range("O10").Value = 0.01588628762541806
$ws.Range("Q10").Value = 0.1714046822742475
$ws.Range("R10").Value = 0.0794314381270903
$ws.Range("S10").Value = 0.387123745819398
$ws.Range("G11").Value = 0.1287425149700599
$ws.Range("J11").Value = 0.09580838323353294
$ws.Range("K11").Value = 0.1976047904191617
$ws.Range("L11").Value = 0.5508982035928144
$ws.Range("S11").Value = 0.02694610778443114
$ws.Range("G12").Value = 0.6963350785340314
$ws.Range("J12").Value = 0.2303664921465969
$ws.Range("L12").Value = 0.03664921465968586
$ws.Range("S12").Value = 0.03664921465968586
$ws.Range("G13").Value = 0.6206896551724138
$ws.Range("J13").Value = 0.2758620689655172
$ws.Range("S13").Value = 0.103448275862069
$ws.Range("F15").Value = 0.01904761904761905
$ws.Range("H15").Value = 0.1571428571428571
$ws.Range("I15").Value = 0.04285714285714286
$ws.Range("J15").Value = 0.3428571428571429
$ws.Range("K15").Value = 0.1
$ws.Range("M15").Value = 0.01428571428571429
$ws.Range("N15").Value = 0.004761904761904762
$ws.Range("O15").Value = 0.05714285714285714
$ws.Range("S15").Value = 0.2619047619047619
$ws.Range("F16").Value = 0.01408450704225352
$ws.Range("H16").Value = 0.2065727699530517
$ws.Range("I16").Value = 0.107981220657277
$ws.Range("J16").Value = 0.3145539906103286
$ws.Range("K16").Value = 0.1126760563380282
$ws.Range("M16").Value = 0.04694835680751173
$ws.Range("N16").Value = 0.004694835680751174
$ws.Range("O16").Value = 0.04694835680751173
$ws.Range("S16").Value = 0.1455399061032864
$ws.Range("F17").Value = 0.01047120418848168
$ws.Range("H17").Value = 0.1649214659685864
$ws.Range("I17").Value = 0.07329842931937172
$ws.Range("J17").Value = 0.3926701570680629
$ws.Range("K17").Value = 0.09424083769633508
$ws.Range("M17").Value = 0.03141361256544502
$ws.Range("N17").Value = 0.002617801047120419
$ws.Range("O17").Value = 0.07068062827225131
$ws.Range("S17").Value = 0.1596858638743456
$ws.Range("F18").Value = 0.03589743589743589
$ws.Range("H18").Value = 0.1384615384615385
$ws.Range("I18").Value = 0.05641025641025641
$ws.Range("J18").Value = 0.4
$ws.Range("K18").Value = 0.1384615384615385
$ws.Range("M18").Value = 0.01025641025641026
$ws.Range("N18").Value = 0.005128205128205128
$ws.Range("O18").Value = 0.05641025641025641
$ws.Range("S18").Value = 0.158974358974359
$ws.Range("F19").Value = 0.01160261058738216
$ws.Range("H19").Value = 0.2422044960116026
$ws.Range("I19").Value = 0.07469180565627266
$ws.Range("J19").Value = 0.321972443799855
$ws.Range("K19").Value = 0.1131254532269761
$ws.Range("M19").Value = 0.02393038433647571
$ws.Range("N19").Value = 0.0007251631617113851
$ws.Range("O19").Value = 0.06526468455402465
$ws.Range("S19").Value = 0.1464829586656998
